# "minor edits to the plan"
#
# Five small changes to Eager/paper/socc15/EditPlan.docx:
#   1) drop the stray _GoBack bookmark after the "Shephard: ..." line
#   2) de-duplicate "paper paper" -> "paper" in the QBETS-contribution sentence
#   3) "quantile" -> "percentile" in the QBETS measurement-set paragraph
#   4) "that were well documented" -> "that was well documented" (GitHub apps)
#   5) "10000"/"1000000" -> "10,000"/"1000,000" and the _GoBack bookmark
#      reappears there (last edited spot)

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that sits right after the
#    "Shephard: ..." paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) "...the contribution that we make with this paper paper to cloud
#    computing." -> drop the duplicated "paper".
# ------------------------------------------------------------------
$d.Content.Find.Execute("paper paper to cloud computing", $true, $true, $false, $false, $false, $true, 1, $false, "paper to cloud computing", 2)

# ------------------------------------------------------------------
# 3) "...confidence bound specified and the quantile of interest..."
#    -> "percentile"
# ------------------------------------------------------------------
$d.Content.Find.Execute("quantile", $true, $true, $false, $false, $false, $true, 1, $false, "percentile", 2)

# ------------------------------------------------------------------
# 4) "...picked a set that were well documented, builds and runs
#    without errors..." -> "that was"
# ------------------------------------------------------------------
$d.Content.Find.Execute("picked a set that were well documented", $true, $true, $false, $false, $false, $true, 1, $false, "picked a set that was well documented", 2)

# ------------------------------------------------------------------
# 5) "...if necessary (e.g. 10000 or 1000000 entities)..." -> add
#    thousands separators, and re-insert the _GoBack bookmark right
#    after the first comma (mirrors where Word would have left it
#    after the last edit made to the document).
# ------------------------------------------------------------------
$d.Content.Find.Execute("10000 or 1000000 entities", $true, $true, $false, $false, $false, $true, 1, $false, "10,000 or 1000,000 entities", 2)

$full = $d.Content
$found = $full.Find.Execute("10,", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $bmRange = $d.Range($full.End, $full.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
